$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '27.456.76'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '1.574.93'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  -0.22%  '
Set-TextValue 'D5' '207.43'
$ws.Range('E5').Value = '  +0.34%  '
Set-TextValue 'D6' '0.499'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  -0.19%  '
Set-TextValue 'D8' '22.30'
$ws.Range('E8').Value = '  +0.73%  '
Set-TextValue 'D9' '0.250'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '1.797.16'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').Value = '1.579.02'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('E14').Value = '  -0.59%  '
Set-TextValue 'D15' '0.521'
$ws.Range('E15').Value = '  -1.60%  '
Set-TextValue 'D16' '63.53'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '27.452.09'
$ws.Range('E17').Value = '  -0.65%  '
Set-TextValue 'D18' '214.14'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('E21').Value = '  -0.23%  '
Set-TextValue 'D22' '4.15'
$ws.Range('E22').Value = '  +0.60%  '
Set-TextValue 'D23' '9.51'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('E24').Value = '  +1.16%  '
Set-TextValue 'D25' '153.41'
$ws.Range('E25').Value = '  -0.34%  '
Set-TextValue 'D26' '6.69'
$ws.Range('E26').Value = '  -2.57%  '
Set-TextValue 'D27' '14.98'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('D33').Value = '1.398.41'
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('E34').Value = '  +1.07%  '
$ws.Range('E35').Value = '  +2.44%  '
Set-TextValue 'D36' '2.30'
$ws.Range('E36').Value = '  -0.25%  '
Set-TextValue 'D37' '0.939'
$ws.Range('E37').Value = '  -3.37%  '
$ws.Range('E38').Value = '  -0.29%  '
Set-TextValue 'D39' '0.531'
$ws.Range('E39').Value = '  -0.68%  '
Set-TextValue 'D40' '0.826'
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('E41').Value = '  -0.20%  '
Set-TextValue 'D42' '1.00'
$ws.Range('E42').Value = '  +2.68%  '
$ws.Range('E43').Value = '  +5.17%  '
Set-TextValue 'D44' '64.55'
$ws.Range('E44').Value = '  +1.28%  '
Set-TextValue 'D45' '2.19'
$ws.Range('E45').Value = '  +0.57%  '
Set-TextValue 'D46' '5.26'
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('D47').Value = '1.709.50'
$ws.Range('E47').Value = '  -0.61%  '
Set-TextValue 'D48' '86.00'
$ws.Range('E48').Value = '  -2.43%  '
$ws.Range('D49').Value = '0.0₇0998'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('E51').Value = '  -0.68%  '
